$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.371.31'
$ws.Range("E2").Value = '  +0.04%  '

$ws.Range("D3").Value = '2.622.30'
$ws.Range("E3").Value = '  -1.71%  '

$ws.Range("E4").Value = '  +0.05%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '593.73'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.74%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '167.90'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.24%  '

$ws.Range("E7").Value = '  +0.05%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.533'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.36%  '

$ws.Range("D9").Value = '2.621.79'
$ws.Range("E9").Value = '  -1.68%  '

$ws.Range("E10").Value = '  -1.96%  '

$ws.Range("E11").Value = '  +1.13%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.364'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.78%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.22'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.05%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '27.61'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.47%  '

$ws.Range("D15").Value = '3.104.37'
$ws.Range("E15").Value = '  -1.52%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.0000182'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.93%  '

$ws.Range("D17").Value = '67.281.59'
$ws.Range("E17").Value = '  +0.15%  '

$ws.Range("D18").Value = '2.621.42'
$ws.Range("E18").Value = '  -1.84%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.01'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.63%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '8.00'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +4.67%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '355.85'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.99%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.30'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.41%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.65'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.91%  '

$ws.Range("E24").Value = '  -0.04%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.92'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -4.92%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '10.25'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +2.24%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '69.56'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.48%  '

$ws.Range("D28").Value = '2.760.45'
$ws.Range("E28").Value = '  -1.40%  '

$ws.Range("E29").Value = '  +0.06%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0000100'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.73%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '542.90'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.20%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '7.89'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.32%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.34'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -3.10%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.89'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.68%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.135'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +4.41%  '

$ws.Range("E36").Value = '  +0.07%  '

$ws.Range("E37").Value = '  -3.37%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '156.38'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.16%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '18.99'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.73%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.365'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.92%  '

$ws.Range("E41").Value = '  -1.06%  '

$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '18.20'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.40%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.20'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.87%  '

$ws.Range("E44").Value = '  +0.06%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.83%  '

$ws.Range("D46").Value = '0.0₆0296'
$ws.Range("E46").Value = '  -0.39%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '152.47'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.23%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.579'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.72%  '

$ws.Range("E49").Value = '  -1.42%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.69'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.25%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0768'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.44%  '
